# Generate Report for Handoff
# -----------------------------------------------------------------------
# The "4f3491b5..." file moved from "Handed back" to "Ready for handoff"
# (with a fresh handoff timestamp), and the report rows for the two files
# were re-sorted (the still-"Handed back" a53b297b file now sorts first,
# row 2; the newly re-handed-off 4f3491b5 file drops to row 3) on every
# sheet. This script rewrites each sheet's data rows + matching hyperlink
# display text to reflect that.

$wb = $excel.ActiveWorkbook

$fileA = "a53b297b-d811-421c-9e68-f34339466385.md"
$fileB = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"

$handedBack = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row 2 -> a53b297b (still handed back),
#                 row 3 -> 4f3491b5 (now ready for handoff, new date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $fileA
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsOverview.Range("D2").Value = "2016-03-23 20:53:40"

$wsOverview.Range("A3").Value = $fileB
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff
$wsOverview.Range("D3").Value = "2016-03-23 20:55:11"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $fileA
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $fileB
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhAFile = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf"
$zhBFile = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf"

# Row 2 -> a53b297b, still "Handed back", handoff datetime unchanged
$wsZh.Range("A2").Value = $fileA
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("D2").Value = $zhAFile
$wsZh.Range("E2").Value = "2016-03-23 20:53:35"
$wsZh.Range("F2").Value = $fileA
$wsZh.Range("G2").Value = $zhAFile
$wsZh.Range("H2").Value = "2016-03-23 20:54:13"
$wsZh.Range("J2").Value = "Include"

# Row 3 -> 4f3491b5, now "Ready for handoff" with a new handoff datetime
$wsZh.Range("A3").Value = $fileB
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = $zhBFile
$wsZh.Range("E3").Value = "2016-03-23 20:55:06"
$wsZh.Range("F3").Value = $fileB
$wsZh.Range("G3").Value = $zhBFile
$wsZh.Range("H3").Value = "2016-03-23 20:54:13"
$wsZh.Range("J3").Value = "Include"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $fileA
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = $zhAFile
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = $fileA
    } elseif ($addr -eq '$G$2') {
        $h.TextToDisplay = $zhAFile
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $fileB
    } elseif ($addr -eq '$D$3') {
        $h.TextToDisplay = $zhBFile
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = $fileB
    } elseif ($addr -eq '$G$3') {
        $h.TextToDisplay = $zhBFile
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deAFile = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf"
$deBFile = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf"

# Row 2 -> a53b297b, still "Handed back", handoff datetime unchanged
$wsDe.Range("A2").Value = $fileA
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("D2").Value = $deAFile
$wsDe.Range("E2").Value = "2016-03-23 20:53:40"
$wsDe.Range("F2").Value = $fileA
$wsDe.Range("G2").Value = $deAFile
$wsDe.Range("H2").Value = "2016-03-23 20:54:22"
$wsDe.Range("J2").Value = "Include"

# Row 3 -> 4f3491b5, now "Ready for handoff" with a new handoff datetime
$wsDe.Range("A3").Value = $fileB
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = $deBFile
$wsDe.Range("E3").Value = "2016-03-23 20:55:11"
$wsDe.Range("F3").Value = $fileB
$wsDe.Range("G3").Value = $deBFile
$wsDe.Range("H3").Value = "2016-03-23 20:54:22"
$wsDe.Range("J3").Value = "Include"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $fileA
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = $deAFile
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = $fileA
    } elseif ($addr -eq '$G$2') {
        $h.TextToDisplay = $deAFile
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $fileB
    } elseif ($addr -eq '$D$3') {
        $h.TextToDisplay = $deBFile
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = $fileB
    } elseif ($addr -eq '$G$3') {
        $h.TextToDisplay = $deBFile
    }
}

Write-Output "done"
